# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12411
$ws1.Range("F3").Value = 588
$ws1.Range("F4").Value = 2017
$ws1.Range("F5").Value = 257
$ws1.Range("F7").Value = 232
$ws1.Range("F8").Value = 12348
$ws1.Range("F9").Value = 2423
$ws1.Range("F11").Value = 9
$ws1.Range("F13").Value = 14
$ws1.Range("F15").Value = 635
$ws1.Range("F16").Value = 2823
$ws1.Range("F17").Value = 6037
$ws1.Range("F19").Value = 3589
$ws1.Range("F20").Value = 211

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12411
$ws4.Range("F3").Value = 588
$ws4.Range("F4").Value = 2018
$ws4.Range("F5").Value = 257
$ws4.Range("F8").Value = 232
$ws4.Range("F9").Value = 12348
$ws4.Range("F10").Value = 2425
$ws4.Range("F12").Value = 9
$ws4.Range("F14").Value = 14
$ws4.Range("F16").Value = 635
$ws4.Range("F17").Value = 2823
$ws4.Range("F19").Value = 6037
$ws4.Range("F21").Value = 3589
$ws4.Range("F22").Value = 211
